$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update progression percentages ---
$ws.Range("E3").Value = 0.4
$ws.Range("E5").Value = 0.35

# --- Insert a new row 7 for the new task "Remettre le code aux normes" ---
# Insert a blank row at position 7 (shifts existing rows 7-16 down to 8-17)
$ws.Rows(7).Insert()

# Copy the formatting of row 6 (the row just above) onto the newly inserted row 7,
# so it keeps the same cell styles used elsewhere in the sheet (no new styles created)
$ws.Range("A6:F6").Copy()
$ws.Range("A7:F7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the content of the new row 7
$ws.Range("A7").Value = "Remettre le code aux normes"
$ws.Range("B7").Value = "Respecter les règles de nommage d'après le document"
$ws.Range("C7").Value = "Jo / Pizzi "
$ws.Range("D7").Value = "En attente "
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = "normal"

# --- Update the sheet view: clear the scrolled top-left cell and move the selection to E5 ---
[void]$ws.Range("E5").Select()
